# Workbook/worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Strip the bespoke header/border/fill formatting that used to live in
#    styles.xml (bold header row, bordered/wrap H column, purple highlight
#    fills). Resetting the whole used block (A1:I40) back to the "Normal"
#    cell style removes every explicit style index from the cells, which is
#    what the target file does.
# ---------------------------------------------------------------------------
$ws.Range("A1:I40").Style = "Normal"

# Column I was only ever used to carry the header's style - with the style
# gone the column has no real content left, so clear it outright (this also
# shrinks the sheet dimension from A1:I.. down to A1:H..).
$ws.Columns.Item(9).Clear()

# Remove the custom row heights (ht="16") left over from the old bold/boxed
# header styling; AutoFit drops the explicit height so rows fall back to the
# sheet's default row height.
$ws.Range("A1:H41").EntireRow.AutoFit()

# ---------------------------------------------------------------------------
# 2) Append the new fastq/metadata row (row 41) that the commit adds.
# ---------------------------------------------------------------------------
# Force text storage for the date-like / text columns so "10.18.18" etc.
# are written as shared strings instead of being auto-parsed into serial
# date numbers.
$ws.Range("A41").NumberFormat = "@"
$ws.Range("B41").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("H41").NumberFormat = "@"

$ws.Range("A41").Value2 = "10.18.18"
$ws.Range("B41").Value2 = "H.BROWN"
$ws.Range("C41").Value2 = 41
$ws.Range("D41").Value2 = "10.18.18"
$ws.Range("E41").Value2 = "H.BROWN"
$ws.Range("F41").Value2 = 41
$ws.Range("G41").Value2 = "ATGACAG"
$ws.Range("H41").Value2 = "E7420L"

# Keep the new row free of any explicit style too.
$ws.Range("A41:H41").Style = "Normal"
$ws.Range("A41:H41").EntireRow.AutoFit()

# ---------------------------------------------------------------------------
# 3) Move the live selection the way the author left it.
# ---------------------------------------------------------------------------
$ws.Range("I50").Select()
